$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# Copy the formatting from the last existing row (row 4) down to the new
# row (row 5) so the new cells inherit the same styles (date format on
# column A, centered text on the rest) without creating new style entries.
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new task row (row 5): Date, Task Name, Status, Person
$ws.Cells.Item(5, 1).Value = 45616
$ws.Cells.Item(5, 2).Value = "PA2"
$ws.Cells.Item(5, 3).Value = "Done"
$ws.Cells.Item(5, 4).Value = "Drew Hutchinson"

$ws.Range("E5").Select()

$wb.Save()
